$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (cell, newValue) updates
# Generated from the authoritative diff for commit 'Add data for 2025-10-23'

$updates = @{}

$updates['Citywide Totals'] = @{ 'L2' = 5453; 'L3' = 5933; 'C4' = 1872; 'E4' = 2060; 'G4' = 1509; 'K4' = 1794; 'L4' = 1448; 'L5' = 357; 'L6' = 4881; 'C7' = 28416; 'E7' = 26065; 'G7' = 24736; 'K7' = 27587; 'L7' = 18072 }
$updates['Logan Square'] = @{ 'L6' = 67; 'L7' = 200 }
$updates['Austin'] = @{ 'L3' = 421; 'L5' = 41; 'L7' = 1197 }
$updates['South Chicago'] = @{ 'L3' = 161; 'L7' = 400 }
$updates['Garfield Park'] = @{ 'L3' = 294; 'G4' = 51; 'L4' = 54; 'G7' = 1428 }
$updates['Woodlawn'] = @{ 'L6' = 70; 'L7' = 313 }
$updates['By Neighborhood'] = @{ 'L6' = 138; 'L7' = 591; 'L8' = 1197; 'L9' = 105; 'E11' = 371; 'L11' = 294; 'L20' = 448; 'L23' = 198; 'L24' = 47; 'L25' = 110; 'L27' = 157; 'L29' = 1019; 'G33' = 1428; 'L36' = 233; 'C42' = 1148; 'L42' = 586; 'L48' = 232; 'L51' = 229; 'L52' = 367; 'L53' = 200; 'L54' = 391; 'L55' = 191; 'L59' = 31; 'E63' = 392; 'K63' = 179; 'L63' = 58; 'L67' = 624; 'L70' = 30; 'L72' = 72; 'L73' = 144; 'L79' = 492; 'L83' = 400; 'L90' = 189; 'L91' = 237; 'L94' = 221; 'L99' = 313; 'C101' = 28416; 'E101' = 26065; 'G101' = 24736; 'K101' = 27587; 'L101' = 18072 }
$updates['North Lawndale'] = @{ 'L3' = 240; 'L7' = 624 }
$updates['Loop'] = @{ 'L6' = 189; 'L7' = 391 }
$updates['Englewood'] = @{ 'L2' = 304; 'L7' = 1019 }
$updates['Lake View'] = @{ 'L6' = 96; 'L7' = 232 }
$updates['Ashburn'] = @{ 'L2' = 57; 'L7' = 138 }
$updates['Humboldt Park'] = @{ 'L3' = 203; 'C4' = 57; 'L6' = 159; 'C7' = 1148; 'L7' = 586 }
$updates['Lower West Side'] = @{ 'L3' = 65; 'L7' = 191 }
$updates['Dunning'] = @{ 'L2' = 19; 'L7' = 47 }
$updates['Douglas'] = @{ 'L3' = 77; 'L7' = 198 }
$updates['Washington Park'] = @{ 'L3' = 108; 'L7' = 237 }
$updates['Roseland'] = @{ 'L6' = 125; 'L7' = 492 }
$updates['Chicago Lawn'] = @{ 'L2' = 139; 'L3' = 151; 'L7' = 448 }
$updates['Grand Boulevard'] = @{ 'L2' = 82; 'L7' = 233 }
$updates['Auburn Gresham'] = @{ 'L2' = 202; 'L5' = 16; 'L7' = 591 }
$updates['West Loop'] = @{ 'L6' = 86; 'L7' = 221 }
$updates['East Side'] = @{ 'L3' = 49; 'L7' = 110 }
$updates['Belmont Cragin'] = @{ 'E4' = 24; 'L4' = 20; 'E7' = 371; 'L7' = 294 }
$updates['Avalon Park'] = @{ 'L3' = 42; 'L7' = 105 }
$updates['Portage Park'] = @{ 'L2' = 52; 'L7' = 144 }
$updates['Montclare'] = @{ 'L3' = 13; 'L7' = 31 }
$updates['West Town'] = @{ 'L2' = 36; 'L6' = 69 }
$updates['O''Hare'] = @{ 'L2' = 13; 'L7' = 30 }
$updates['Edgewater'] = @{ 'L2' = 44; 'L7' = 157 }
$updates['Washington Heights'] = @{ 'L2' = 63; 'L6' = 52; 'L7' = 189 }
$updates['Little Italy, UIC'] = @{ 'L3' = 71; 'L7' = 229 }
$updates['Old Town'] = @{ 'L3' = 19; 'L7' = 72 }
$updates['Little Village'] = @{ 'L2' = 119; 'L7' = 367 }

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}

Write-Output "Applied updates to $($updates.Keys.Count) sheets"
